$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; remove protection so the cells below can be
# written, then reinstate protection afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidentiality disclaimer (A10)
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-20 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) columns for rows 2-7
$ws.Range("D2").Value = 0.4782961621081183
$ws.Range("E2").Value = 0.006291781360597604

$ws.Range("D3").Value = 0.3403815123887585
$ws.Range("E3").Value = 0.005482041587901776

$ws.Range("D4").Value = 0.09612795612706521
$ws.Range("E4").Value = 0.01135339700846982

$ws.Range("D5").Value = 0.05358521664241381
$ws.Range("E5").Value = 0.004272024015702591

$ws.Range("D6").Value = 0.03160915273364417
$ws.Range("E6").Value = 0.01291837933059292

$ws.Range("E7").Value = 0.006603955691503627

# Restore sheet protection
$ws.Protect()
